# Adds a new "sleep" column (M) with per-participant values, widens
# column L slightly, and moves the active selection to Q7 — matching the
# "enhance participant mapping and additional variables handling" edit.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New header in M1
$ws.Cells.Item(1, 13).Value = "sleep"

# New "sleep" values for each participant row (2-21)
$sleepValues = @(1, 1, 1, 2, 1, 1, 1, 2, 1, 2, 1, 1, 1, 2, 1, 1, 1, 2, 1, 2)

for ($i = 0; $i -lt $sleepValues.Length; $i++) {
    $row = $i + 2
    $ws.Cells.Item($row, 13).Value = $sleepValues[$i]
}

# Widen column L (completion_date) as seen in the authored workbook
$ws.Columns.Item(12).ColumnWidth = 19

# Move the selection to match the saved UI state
$ws.Range("Q7").Select()
